# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new blank column before the
# existing "Late" column (column N), pushing "Late" -> O, the second
# "heading" column -> P, and "Outstanding" -> Q. Also switch the active
# sheet/tab from "NewLoanInput" to "Repayment schedule".

$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (14); existing N:P data shifts to O:Q.
[void]$wsRepay.Columns.Item(14).Insert()

# Excel copies the width of the column to the left (M) onto the freshly
# inserted column, but without the bestFit flag - match that explicitly.
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab with a fresh selection
# (this also clears the previous active-tab flag on "NewLoanInput").
[void]$wsRepay.Activate()
[void]$wsRepay.Range("K17").Select()
